$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Row=42; Col=8; Val=2002.6666},
    @{Row=42; Col=10; Val=0},
    @{Row=42; Col=12; Val=0},
    @{Row=42; Col=14; Val=$null},
    @{Row=70; Col=8; Val=1020},
    @{Row=70; Col=9; Val=778.1111},
    @{Row=70; Col=10; Val=1382.8334},
    @{Row=70; Col=11; Val=2334.3333},
    @{Row=70; Col=12; Val=4148.5002},
    @{Row=70; Col=13; Val=-2064.3333},
    @{Row=70; Col=14; Val=-4688.5002},
    @{Row=73; Col=8; Val=1020},
    @{Row=73; Col=9; Val=778.1111},
    @{Row=73; Col=10; Val=1382.8334},
    @{Row=73; Col=11; Val=2334.3333},
    @{Row=73; Col=12; Val=4148.5002},
    @{Row=73; Col=13; Val=-1398.3333},
    @{Row=73; Col=14; Val=-6020.5002},
    @{Row=76; Col=8; Val=142866530},
    @{Row=76; Col=9; Val=200008140},
    @{Row=76; Col=10; Val=12500},
    @{Row=76; Col=11; Val=200008140},
    @{Row=76; Col=12; Val=12500},
    @{Row=76; Col=13; Val=-200007825},
    @{Row=76; Col=14; Val=-13130},
    @{Row=79; Col=8; Val=142866530},
    @{Row=79; Col=9; Val=200008140},
    @{Row=79; Col=10; Val=12500},
    @{Row=79; Col=11; Val=200008140},
    @{Row=79; Col=12; Val=12500},
    @{Row=79; Col=13; Val=-200007048},
    @{Row=79; Col=14; Val=-14684},
    @{Row=80; Col=8; Val=2258},
    @{Row=80; Col=9; Val=201.25},
    @{Row=80; Col=10; Val=4314.75},
    @{Row=80; Col=11; Val=603.75},
    @{Row=80; Col=12; Val=12944.25},
    @{Row=80; Col=13; Val=394.25},
    @{Row=80; Col=14; Val=-14940.25},
    @{Row=83; Col=8; Val=2258},
    @{Row=83; Col=9; Val=201.25},
    @{Row=83; Col=10; Val=4314.75},
    @{Row=83; Col=11; Val=1811.25},
    @{Row=83; Col=12; Val=38832.75},
    @{Row=83; Col=13; Val=3180.75},
    @{Row=83; Col=14; Val=-48816.75},
    @{Row=99; Col=8; Val=90912260},
    @{Row=99; Col=9; Val=808.1111},
    @{Row=99; Col=10; Val=500013760},
    @{Row=99; Col=11; Val=2424.3333},
    @{Row=99; Col=12; Val=1500041280},
    @{Row=99; Col=13; Val=-926.3332999999998},
    @{Row=99; Col=14; Val=-1500044276},
    @{Row=132; Col=8; Val=265220.5},
    @{Row=132; Col=9; Val=310571.97},
    @{Row=132; Col=11; Val=931715.9099999999},
    @{Row=132; Col=13; Val=-929185.9099999999},
    @{Row=137; Col=8; Val=7036.278},
    @{Row=137; Col=9; Val=4180.8213},
    @{Row=137; Col=10; Val=17030.375},
    @{Row=137; Col=11; Val=12542.4639},
    @{Row=137; Col=12; Val=51091.125},
    @{Row=137; Col=13; Val=-9992.463899999999},
    @{Row=137; Col=14; Val=-56191.125},
    @{Row=138; Col=8; Val=3577.5925},
    @{Row=138; Col=9; Val=1054.0714},
    @{Row=138; Col=10; Val=6295.231},
    @{Row=138; Col=11; Val=3162.2142},
    @{Row=138; Col=12; Val=18885.693},
    @{Row=138; Col=13; Val=1977.7858},
    @{Row=138; Col=14; Val=-29165.693},
    @{Row=140; Col=8; Val=219999},
    @{Row=140; Col=10; Val=333332},
    @{Row=140; Col=12; Val=333332},
    @{Row=140; Col=14; Val=-343692},
    @{Row=141; Col=8; Val=873.1111},
    @{Row=141; Col=9; Val=790.119},
    @{Row=141; Col=11; Val=2370.357},
    @{Row=141; Col=13; Val=2809.643}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Row=32; Col=8; Val=13893531},
    @{Row=32; Col=9; Val=13892719},
    @{Row=32; Col=10; Val=13897184},
    @{Row=32; Col=11; Val=13892719},
    @{Row=32; Col=12; Val=13897184},
    @{Row=32; Col=13; Val=-13892432},
    @{Row=32; Col=14; Val=-13897758},
    @{Row=61; Col=8; Val=3114.6943},
    @{Row=61; Col=9; Val=3061.1143},
    @{Row=61; Col=10; Val=4990},
    @{Row=61; Col=11; Val=3061.1143},
    @{Row=61; Col=12; Val=4990},
    @{Row=61; Col=13; Val=-2849.1143},
    @{Row=61; Col=14; Val=-5414},
    @{Row=63; Col=8; Val=1492.6},
    @{Row=63; Col=9; Val=921},
    @{Row=63; Col=10; Val=2350},
    @{Row=63; Col=11; Val=921},
    @{Row=63; Col=12; Val=2350},
    @{Row=63; Col=13; Val=-235},
    @{Row=63; Col=14; Val=-3722},
    @{Row=66; Col=8; Val=1492.6},
    @{Row=66; Col=9; Val=921},
    @{Row=66; Col=10; Val=2350},
    @{Row=66; Col=11; Val=4605},
    @{Row=66; Col=12; Val=11750},
    @{Row=66; Col=13; Val=-1173},
    @{Row=66; Col=14; Val=-18614},
    @{Row=74; Col=8; Val=5048.15},
    @{Row=74; Col=9; Val=6409.727},
    @{Row=74; Col=11; Val=6409.727},
    @{Row=74; Col=13; Val=-5535.727},
    @{Row=77; Col=8; Val=5048.15},
    @{Row=77; Col=9; Val=6409.727},
    @{Row=77; Col=11; Val=32048.635},
    @{Row=77; Col=13; Val=-27680.635},
    @{Row=97; Col=8; Val=922.6774},
    @{Row=97; Col=9; Val=624.381},
    @{Row=97; Col=10; Val=1549.1},
    @{Row=97; Col=11; Val=624.381},
    @{Row=97; Col=12; Val=1549.1},
    @{Row=97; Col=13; Val=-128.381},
    @{Row=97; Col=14; Val=-2541.1},
    @{Row=102; Col=8; Val=2100},
    @{Row=102; Col=9; Val=1562.6364},
    @{Row=102; Col=11; Val=1562.6364},
    @{Row=102; Col=13; Val=59.36359999999991},
    @{Row=122; Col=8; Val=4162.0586},
    @{Row=122; Col=9; Val=3765},
    @{Row=122; Col=11; Val=11295},
    @{Row=122; Col=13; Val=-8845},
    @{Row=132; Col=8; Val=342051.1},
    @{Row=132; Col=9; Val=377987},
    @{Row=132; Col=11; Val=1133961},
    @{Row=132; Col=13; Val=-1131431},
    @{Row=136; Col=8; Val=3114.6943},
    @{Row=136; Col=9; Val=3061.1143},
    @{Row=136; Col=10; Val=4990},
    @{Row=136; Col=11; Val=9183.3429},
    @{Row=136; Col=12; Val=14970},
    @{Row=136; Col=13; Val=-6633.3429},
    @{Row=136; Col=14; Val=-20070}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Row=62; Col=8; Val=21300},
    @{Row=62; Col=9; Val=19833.334},
    @{Row=62; Col=10; Val=23500},
    @{Row=62; Col=11; Val=19833.334},
    @{Row=62; Col=12; Val=23500},
    @{Row=62; Col=13; Val=-19209.334},
    @{Row=62; Col=14; Val=-24748},
    @{Row=65; Col=8; Val=21300},
    @{Row=65; Col=9; Val=19833.334},
    @{Row=65; Col=10; Val=23500},
    @{Row=65; Col=11; Val=99166.67},
    @{Row=65; Col=12; Val=117500},
    @{Row=65; Col=13; Val=-96046.67},
    @{Row=65; Col=14; Val=-123740},
    @{Row=105; Col=8; Val=46824924},
    @{Row=105; Col=9; Val=59831092},
    @{Row=105; Col=11; Val=59831092},
    @{Row=105; Col=13; Val=-59829345},
    @{Row=134; Col=8; Val=43488416},
    @{Row=134; Col=10; Val=22106.857},
    @{Row=134; Col=12; Val=66320.571},
    @{Row=134; Col=14; Val=-71390.571}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Row=121; Col=8; Val=18224},
    @{Row=121; Col=10; Val=25989},
    @{Row=121; Col=12; Val=77967},
    @{Row=121; Col=14; Val=-80587}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Row=122; Col=8; Val=49285.953},
    @{Row=122; Col=9; Val=127831.25},
    @{Row=122; Col=11; Val=383493.75},
    @{Row=122; Col=13; Val=-381043.75},
    @{Row=132; Col=8; Val=32262064},
    @{Row=132; Col=9; Val=40003268},
    @{Row=132; Col=11; Val=120009804},
    @{Row=132; Col=13; Val=-120007274}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Row=68; Col=8; Val=1327.6666},
    @{Row=68; Col=10; Val=1241},
    @{Row=68; Col=12; Val=1241},
    @{Row=68; Col=14; Val=-2739},
    @{Row=71; Col=8; Val=1327.6666},
    @{Row=71; Col=10; Val=1241},
    @{Row=71; Col=12; Val=6205},
    @{Row=71; Col=14; Val=-13693},
    @{Row=132; Col=8; Val=4081.6667},
    @{Row=132; Col=9; Val=3535.75},
    @{Row=132; Col=11; Val=10607.25},
    @{Row=132; Col=13; Val=-8077.25},
    @{Row=136; Col=8; Val=26318640},
    @{Row=136; Col=9; Val=10872334},
    @{Row=136; Col=10; Val=90912290},
    @{Row=136; Col=11; Val=32617002},
    @{Row=136; Col=12; Val=272736870},
    @{Row=136; Col=13; Val=-32614452},
    @{Row=136; Col=14; Val=-272741970}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Row=62; Col=8; Val=14399.833},
    @{Row=62; Col=9; Val=8729.799999999999},
    @{Row=62; Col=10; Val=42750},
    @{Row=62; Col=11; Val=8729.799999999999},
    @{Row=62; Col=12; Val=42750},
    @{Row=62; Col=13; Val=-8105.799999999999},
    @{Row=62; Col=14; Val=-43998},
    @{Row=65; Col=8; Val=14399.833},
    @{Row=65; Col=9; Val=8729.799999999999},
    @{Row=65; Col=10; Val=42750},
    @{Row=65; Col=11; Val=43649},
    @{Row=65; Col=12; Val=213750},
    @{Row=65; Col=13; Val=-40529},
    @{Row=65; Col=14; Val=-219990},
    @{Row=100; Col=8; Val=1887.5555},
    @{Row=100; Col=9; Val=1051.7333},
    @{Row=100; Col=11; Val=2103.4666},
    @{Row=100; Col=13; Val=-1562.4666},
    @{Row=132; Col=8; Val=3712.4285},
    @{Row=132; Col=9; Val=3054.9539},
    @{Row=132; Col=10; Val=12259.6},
    @{Row=132; Col=11; Val=9164.861699999999},
    @{Row=132; Col=12; Val=36778.8},
    @{Row=132; Col=13; Val=-6634.861699999999},
    @{Row=132; Col=14; Val=-41838.8},
    @{Row=136; Col=8; Val=10009363},
    @{Row=136; Col=9; Val=11633771},
    @{Row=136; Col=11; Val=34901313},
    @{Row=136; Col=13; Val=-34898763},
    @{Row=138; Col=8; Val=129000},
    @{Row=138; Col=10; Val=129000},
    @{Row=138; Col=12; Val=129000},
    @{Row=138; Col=14; Val=-139280},
    @{Row=140; Col=8; Val=88779.60000000001},
    @{Row=140; Col=10; Val=88779.60000000001},
    @{Row=140; Col=12; Val=88779.60000000001},
    @{Row=140; Col=14; Val=-99139.60000000001}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
}
